# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" and populate it with
#    the quarterly fund-holding detail rows.
# 2) Update the "总计" (summary) sheet so the new 2022-Q4 aggregate appears as
#    its first data row, with the remaining quarters shifted down by one row.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) New "2022-Q4" sheet, positioned right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Copy the header style (bold / centred / bordered) from the "总计" sheet's
# header cell onto the new header row.
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Data rows (A = numeric index, H = numeric rank; B/C/D/E/F/G are text, as in
# the other quarter sheets, so leading zeros in fund codes are preserved).
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'007128"
$q4.Range("C2").Value = "天弘增强回报债券A"
$q4.Range("D2").Value = "'44.39"
$q4.Range("E2").Value = "'86.57"
$q4.Range("F2").Value = "'4.56"
$q4.Range("G2").Value = "'2.0242"
$q4.Range("H2").Value = 2

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'007129"
$q4.Range("C3").Value = "天弘增强回报债券C"
$q4.Range("D3").Value = "'43.01"
$q4.Range("E3").Value = "'86.57"
$q4.Range("F3").Value = "'4.56"
$q4.Range("G3").Value = "'1.9613"
$q4.Range("H3").Value = 2

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'009735"
$q4.Range("C4").Value = "天弘增强回报债券E"
$q4.Range("D4").Value = "'3.76"
$q4.Range("E4").Value = "'86.57"
$q4.Range("F4").Value = "'4.56"
$q4.Range("G4").Value = "'0.1715"
$q4.Range("H4").Value = 2

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'009327"
$q4.Range("C5").Value = "东兴兴晟混合A"
$q4.Range("D5").Value = "'0.38"
$q4.Range("E5").Value = "'79.79"
$q4.Range("F5").Value = "'1.09"
$q4.Range("G5").Value = "'0.0041"
$q4.Range("H5").Value = 4

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'009328"
$q4.Range("C6").Value = "东兴兴晟混合C"
$q4.Range("D6").Value = "'0.08"
$q4.Range("E6").Value = "'79.79"
$q4.Range("F6").Value = "'1.09"
$q4.Range("G6").Value = "'0.0009"
$q4.Range("H6").Value = 4

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'002630"
$q4.Range("C7").Value = "江信瑞福灵活配置混合A"
$q4.Range("D7").Value = "'0.01"
$q4.Range("E7").Value = "'86.57"
$q4.Range("F7").Value = "'4.56"
$q4.Range("G7").Value = "'0.0005"
$q4.Range("H7").Value = 2

# Copy the index-column style (bold/centred/bordered) used on column A of the
# other sheets onto the new sheet's A2:A7.
$total.Range("A2").Copy()
$q4.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) "总计" sheet: add the 2022-Q4 row, push the others down one row.
# ---------------------------------------------------------------------
$total.Range("A7").Value = 4
$total.Range("B7").Value = "2021-Q4"
$total.Range("C7").Value = 6
$total.Range("D7").Value = 0.33

$total.Range("A6").Value = 3
$total.Range("B6").Value = "2022-Q1"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.29

$total.Range("A5").Value = 2
$total.Range("B5").Value = "2022-Q2"
$total.Range("C5").Value = 6
$total.Range("D5").Value = 0.29

$total.Range("A4").Value = 1
$total.Range("B4").Value = "2022-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.01

$total.Range("A3").Value = 0
$total.Range("B3").Value = "2022-Q4"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 4.16

Write-Output "2022-Q4 data added"
